$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.665.06'
$ws.Range("E2").Value = '  +1.15%  '
$ws.Range("D3").Value = '1.828.14'
$ws.Range("E3").Value = '  +1.89%  '
$ws.Range("E4").Value = '  +0.36%  '
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = '309.88'
$r.NumberFormat = "General"
$ws.Range("E5").Value = '  +0.92%  '
$ws.Range("E6").Value = '  +0.30%  '
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = '0.4662'
$r.NumberFormat = "General"
$ws.Range("E7").Value = '  +3.38%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  +0.83%  '
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = '0.9048'
$r.NumberFormat = "General"
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = '0.07697'
$r.NumberFormat = "General"
$ws.Range("E11").Value = '  -0.49%  '
$ws.Range("E12").Value = '  -0.11%  '
$ws.Range("D13").Value = '1.857.41'
$ws.Range("E13").Value = '  +3.57%  '
$ws.Range("E14").Value = '  -0.20%  '
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = '6.380'
$r.NumberFormat = "General"
$ws.Range("E15").Value = '  +0.84%  '
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = '87.75'
$r.NumberFormat = "General"
$ws.Range("E16").Value = '  +3.24%  '
$ws.Range("E17").Value = '  +0.26%  '
$ws.Range("E18").Value = '  +0.43%  '
$ws.Range("D20").Value = '26.702.03'
$ws.Range("E20").Value = '  +1.21%  '
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = '14.21'
$r.NumberFormat = "General"
$ws.Range("E21").Value = '  -0.36%  '
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = '5.031'
$r.NumberFormat = "General"
$ws.Range("E22").Value = '  +1.20%  '
$ws.Range("E23").Value = '  +0.13%  '
$ws.Range("E24").Value = '  -2.88%  '
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = '153.01'
$r.NumberFormat = "General"
$ws.Range("E25").Value = '  +1.50%  '
$ws.Range("E26").Value = '  +0.66%  '
$ws.Range("E27").Value = '  -1.23%  '
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = '113.83'
$r.NumberFormat = "General"
$ws.Range("E28").Value = '  +1.64%  '
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = '4.874'
$r.NumberFormat = "General"
$ws.Range("E29").Value = '  +0.56%  '
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = '0.08820'
$r.NumberFormat = "General"
$ws.Range("E30").Value = '  +1.53%  '
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = '3.128'
$r.NumberFormat = "General"
$ws.Range("E31").Value = '  +1.69%  '
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = '2.841'
$r.NumberFormat = "General"
$ws.Range("E32").Value = '  +2.94%  '
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = '1.168'
$r.NumberFormat = "General"
$ws.Range("E33").Value = '  +5.91%  '
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = '0.7358'
$r.NumberFormat = "General"
$ws.Range("E34").Value = '  +1.92%  '
$ws.Range("E35").Value = '  -0.07%  '
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = '1.083'
$r.NumberFormat = "General"
$ws.Range("E36").Value = '  +1.49%  '
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = '0.01934'
$r.NumberFormat = "General"
$ws.Range("E37").Value = '  +0.20%  '
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = '0.05162'
$r.NumberFormat = "General"
$ws.Range("E38").Value = '  +1.19%  '
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = '2.925'
$r.NumberFormat = "General"
$ws.Range("E39").Value = '  +2.31%  '
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = '6.889'
$r.NumberFormat = "General"
$ws.Range("E40").Value = '  +0.94%  '
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = '0.5072'
$r.NumberFormat = "General"
$ws.Range("E41").Value = '  +0.40%  '
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = '0.1500'
$r.NumberFormat = "General"
$ws.Range("E42").Value = '  -1.38%  '
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = '8.081'
$r.NumberFormat = "General"
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = '1.009'
$r.NumberFormat = "General"
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = '0.4666'
$r.NumberFormat = "General"
$ws.Range("E45").Value = '  +0.80%  '
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = '10.06'
$r.NumberFormat = "General"
$ws.Range("E46").Value = '  +2.06%  '
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = '98.56'
$r.NumberFormat = "General"
$ws.Range("E47").Value = '  -2.15%  '
$ws.Range("E48").Value = '  +0.37%  '
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = '0.06034'
$r.NumberFormat = "General"
$ws.Range("E49").Value = '  +1.22%  '
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = '63.97'
$r.NumberFormat = "General"
$ws.Range("E50").Value = '  +0.07%  '
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = '35.79'
$r.NumberFormat = "General"
$ws.Range("E51").Value = '  -0.86%  '
